# "Generate Report for Handoff"
#
# The localization-status report moves from "In Translation" to
# "Ready for handoff" and the associated timestamps are refreshed to the
# handoff-generation time. The per-language detail sheets ("zh-cn",
# "de-de") and the roll-up "Overview" sheet are kept in sync, and the
# Status-ish columns are widened slightly to fit the new, longer text.

$wb = $excel.ActiveWorkbook

# --- "zh-cn" detail sheet: Status + Latest Handoff Datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-16 04:38:01"

# --- "de-de" detail sheet: Status + Latest Handoff Datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-16 04:38:09"

# --- "Overview" roll-up sheet: zh-cn / de-de status + generate date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-16 04:38:09"

# --- Widen the Status-related columns to fit "Ready for handoff" ---
# ColumnWidth is quantized by the host to whole-pixel steps (1/6 of a
# character unit here), so we feed it the pre-image that lands on the
# pixel bucket closest to the desired ~17.22-character width rather than
# the target width itself (which would round to a much wider bucket).
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333336
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333336
